$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, shifting Moorabbin..Springvale rows down by one
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the Melbourne / Nandos exposure site entry
$ws.Cells.Item(16, 1).Value = "Melbourne"
$ws.Cells.Item(16, 2).Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Cells.Item(16, 3).Value = "01/01/2021 2:00am - 2:30am"
$ws.Cells.Item(16, 4).Value = "Case dined at venue"
